# Edit script: replace temperature dataset in "Datos crudos" and refresh
# formulas / view-state so the workbook matches the new 36-row dataset
# (commit: "Mas mediciones de temperatura").

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# 1) Drop the trailing 3 rows (38-40) that no longer exist in the new run.
$ws1.Range("38:40").EntireRow.Delete()

# 2) New timestamp (column C) and temperature (column E) readings for rows 2-37.
$timestamps = @("2023-12-10 07:13:04","2023-12-10 07:14:06","2023-12-10 07:15:08","2023-12-10 07:16:11","2023-12-10 07:17:13","2023-12-10 07:18:16","2023-12-10 07:19:18","2023-12-10 07:20:20","2023-12-10 07:21:23","2023-12-10 07:22:25","2023-12-10 07:23:28","2023-12-10 07:24:30","2023-12-10 07:25:32","2023-12-10 07:26:35","2023-12-10 07:27:37","2023-12-10 07:28:40","2023-12-10 07:29:42","2023-12-10 07:30:44","2023-12-10 07:31:47","2023-12-10 07:32:49","2023-12-10 07:33:52","2023-12-10 07:34:54","2023-12-10 07:35:56","2023-12-10 07:36:59","2023-12-10 07:38:01","2023-12-10 07:39:04","2023-12-10 07:40:06","2023-12-10 07:41:08","2023-12-10 07:42:11","2023-12-10 07:43:13","2023-12-10 07:44:16","2023-12-10 07:45:18","2023-12-10 07:46:20","2023-12-10 07:47:23","2023-12-10 07:48:25","2023-12-10 07:49:28")
$temps = @(29.213286713286699,26.1975524475524,22.6573426573426,21.346153846153801,20.5594405594405,19.9038461538461,18.592657342657301,18.592657342657301,18.199300699300601,17.805944055944,17.019230769230699,16.494755244755201,16.756993006993,16.363636363636299,16.232517482517402,15.9702797202797,15.9702797202797,15.7080419580419,15.7080419580419,15.4458041958042,15.7080419580419,15.839160839160799,15.4458041958042,15.1835664335664,15.4458041958042,15.314685314685301,15.576923076923,15.314685314685301,15.1835664335664,15.4458041958042,15.4458041958042,15.314685314685301,15.576923076923,15.314685314685301,16.756993006993,19.379370629370602)

$n = $timestamps.Length
$cArr = New-Object 'object[,]' $n,1
$eArr = New-Object 'object[,]' $n,1
for ($i = 0; $i -lt $n; $i++) {
    $cArr[$i,0] = $timestamps[$i]
    $eArr[$i,0] = $temps[$i]
}
$ws1.Range("C2:C37").Value = $cArr
$ws1.Range("E2:E37").Value = $eArr

# 3) Offset/count helpers that locate the steady-state window in the new data.
$ws1.Range("H2").Value = 21
$ws1.Range("H3").Formula = "=COUNT(E:E)-2"

# 4) View/selection state: active sheet moves from "Datos validos" to "Datos crudos".
$ws2.Range("B3").Select()
$ws1.Activate()
$ws1.Range("F14").Select()
